$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated 2D training schedule values (rows 2-6, columns B-H)
$data = @{
    2 = @(9, 3, 4, 4, -5, 1, 65)
    3 = @(8, 4, 6, 8, -2, 4, 32)
    4 = @(7, 2, 6, 7, -1, 5, 21)
    5 = @(5, 2, 1, 4, -4, 2, 54)
    6 = @(7, 0, 4, 3, -3, 3, 43)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

# Update the active selection to I1 (no break screen)
$ws.Activate()
$ws.Range("I1").Select()
